$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.808.93'
$ws.Range('E2').Value = '  +2.64%  '
$ws.Range('D3').Value = '3.772.65'
$ws.Range('E3').Value = '  +4.14%  '
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').Value = '427.29'
$ws.Range('E5').Value = '  +6.27%  '
$ws.Range('D6').Value = '138.41'
$ws.Range('E6').Value = '  +10.48%  '
$ws.Range('D7').Value = '0.617'
$ws.Range('E7').Value = '  +4.30%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('D9').Value = '0.734'
$ws.Range('E9').Value = '  +6.13%  '
$ws.Range('D10').Value = '0.151'
$ws.Range('E10').Value = '  -3.53%  '
$ws.Range('D11').Value = '0.0000311'
$ws.Range('E11').Value = '  -8.50%  '
$ws.Range('D12').Value = '42.51'
$ws.Range('E12').Value = '  +7.76%  '
$ws.Range('D13').Value = '10.46'
$ws.Range('E13').Value = '  +11.90%  '
$ws.Range('D14').Value = '4.396.07'
$ws.Range('E14').Value = '  +6.14%  '
$ws.Range('D15').Value = '15.08'
$ws.Range('E15').Value = '  +4.68%  '
$ws.Range('B16').Value = 'TRON'
$ws.Range('C16').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D16').Value = '0.138'
$ws.Range('E16').Value = '  +0.91%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.822.57'
$ws.Range('E17').Value = '  +6.33%  '
$ws.Range('D18').Value = '19.88'
$ws.Range('E18').Value = '  +5.03%  '
$ws.Range('E19').Value = '  +8.44%  '
$ws.Range('D20').Value = '66.107.63'
$ws.Range('E20').Value = '  +3.21%  '
$ws.Range('D21').Value = '405.65'
$ws.Range('E21').Value = '  +1.61%  '
$ws.Range('D22').Value = '15.05'
$ws.Range('E22').Value = '  +6.99%  '
$ws.Range('D23').Value = '3.22'
$ws.Range('E23').Value = '  +8.66%  '
$ws.Range('D24').Value = '84.75'
$ws.Range('E24').Value = '  +2.51%  '
$ws.Range('D25').Value = '36.53'
$ws.Range('E25').Value = '  +5.24%  '
$ws.Range('D26').Value = '3.29'
$ws.Range('E26').Value = '  +7.89%  '
$ws.Range('D27').Value = '9.72'
$ws.Range('E27').Value = '  +40.23%  '
$ws.Range('D28').Value = '9.85'
$ws.Range('E28').Value = '  +10.78%  '
$ws.Range('E29').Value = '  -0.42%  '
$ws.Range('D30').Value = '13.78'
$ws.Range('E30').Value = '  +13.92%  '
$ws.Range('D31').Value = '710.75'
$ws.Range('E31').Value = '  +5.20%  '
$ws.Range('E32').Value = '  +18.05%  '
$ws.Range('E33').Value = '  +4.20%  '
$ws.Range('D34').Value = '40.80'
$ws.Range('E34').Value = '  +9.29%  '
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.15%  '
$ws.Range('D36').Value = '0.149'
$ws.Range('E36').Value = '  -0.85%  '
$ws.Range('E37').Value = '  +36.47%  '
$ws.Range('D38').Value = '56.63'
$ws.Range('E38').Value = '  +4.10%  '
$ws.Range('D39').Value = '0.0474'
$ws.Range('E39').Value = '  +7.28%  '
$ws.Range('D40').Value = '2.64'
$ws.Range('E40').Value = '  +44.48%  '
$ws.Range('E41').Value = '  +8.32%  '
$ws.Range('D42').Value = '0.0₃0681'
$ws.Range('E42').Value = '  +1.73%  '
$ws.Range('E43').Value = '  +6.99%  '
$ws.Range('E44').Value = '  +0.64%  '
$ws.Range('D45').Value = '3.36'
$ws.Range('E45').Value = '  +7.71%  '
$ws.Range('D46').Value = '0.322'
$ws.Range('E46').Value = '  +14.47%  '
$ws.Range('D47').Value = '3.16'
$ws.Range('E47').Value = '  +3.22%  '
$ws.Range('D48').Value = '2.06'
$ws.Range('E48').Value = '  +3.89%  '
$ws.Range('D49').Value = '2.63'
$ws.Range('E49').Value = '  +5.98%  '
$ws.Range('D50').Value = '142.22'
$ws.Range('E50').Value = '  -0.19%  '
$ws.Range('D51').Value = '2.79'
$ws.Range('E51').Value = '  +3.85%  '
